$wb = $excel.ActiveWorkbook

# The handback status report regenerated its timestamps for the
# "Correspond Handoff Datetime" (col E) and "Correspond Handback DateTime"
# (col H) of the first data row (row 2) on both the zh-cn and de-de sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 05:01:46"
$wsZhCn.Range("H2").Value = "2016-03-21 05:02:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 05:01:50"
$wsDeDe.Range("H2").Value = "2016-03-21 05:02:14"
